$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds an Excel date serial value (2023-09-03 => 45172)
# that was bumped forward to 2023-09-06 (45175) for every data row (2..359).
$ws.Range("C2:C359").Value = 45175
